$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Membership counts (column B) are numeric-looking strings stored as text;
# force text format before assigning so Excel does not coerce them to numbers.
$numericTextCells = @("B2","B3","B4","B5","B6")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('B2').Value = '450'
$ws.Range('C2').Value = 'No, FLASCO does not encompass community sites, because it is primarily focused on clinical oncology practices.'
$ws.Range('D2').Value = 'No, FLASCO primarily focuses on education and training initiatives for oncology professionals rather than direct policy advocacy.'
$ws.Range('E2').Value = 'Yes, FLASCO provides engagement opportunity with leadership. FLASCO offers programs and events that allow members to interact with and learn from industry leaders in clinical oncology.'
$ws.Range('F2').Value = 'Yes, FLASCO does provide support for clinical trial recruitment. FLASCO collaborates with research institutions and sponsors to promote clinical trials and facilitate patient access.'
$ws.Range('G2').Value = 'No, FLASCO primarily focuses on education and advocacy for cancer oncologists and clinicians, rather than direct engagement with payors.'
$ws.Range('H2').Value = 'Yes, FLASCO does include area experts on its board. FLASCO is a statewide organization that aims to enhance the quality and effectiveness of oncology care in the state of Florida. This includes having a board comprised of clinical oncologists and other professionals who are considered experts in the field.'
$ws.Range('I2').Value = 'Yes, FLASCO is involved in therapeutic research collaborations. FLASCO actively partners with pharmaceutical companies, academic institutions, and other organizations to conduct research aimed at developing new treatments for cancer patients.'
$ws.Range('J2').Value = 'Yes, FLASCO includes top therapeutic area experts on its board. FLASCO is known to have leading experts in the field of clinical oncology serving on its board, providing valuable insights and guidance for the organization.'
$ws.Range('K2').Value = 'Florida'
$ws.Range('B3').Value = '120'
$ws.Range('C3').Value = 'No, GASCO does not encompass community sites. GASCO focuses on academic and research-based oncology practices, rather than community-based practices.'
$ws.Range('D3').Value = 'No, GASCO is not influential on state or local policy. GASCO is a professional organization focused primarily on education, advocacy, and support for clinical oncology professionals in Georgia. They do not have a direct influence on policy-making decisions at the state or local level.'
$ws.Range('E3').Value = 'Yes, GASCO provides engagement opportunity with leadership. GASCO offers various leadership development programs, mentorship opportunities, and networking events for its members to engage with leaders in the field of oncology.'
$ws.Range('F3').Value = 'No, GASCO does not provide support for clinical trial recruitment. GASCO is a professional organization for oncologists and does not specifically focus on clinical trial recruitment.'
$ws.Range('G3').Value = 'Yes, GASCO provides engagement opportunities with payors. GASCO works closely with payors to ensure that patients have access to high-quality cancer care and treatments.'
$ws.Range('H3').Value = 'No, GASCO does not include area experts on its board. The board members are primarily oncologists and healthcare administrators, rather than specialists in specific areas of oncology.'
$ws.Range('I3').Value = 'Yes, GASCO is involved in therapeutic research collaborations. GASCO actively partners with pharmaceutical companies, research institutions, and healthcare providers to advance cancer treatment options and improve patient outcomes.'
$ws.Range('J3').Value = 'No, GASCO does not include top therapeutic area experts on its board. The primary focus of GASCO is on promoting education and advocacy for clinical oncologists in Georgia, rather than including experts in specific therapeutic areas.'
$ws.Range('K3').Value = 'Southeastern'
$ws.Range('B4').Value = '100'
$ws.Range('C4').Value = 'No, IOS focuses on professional oncology practices, not community sites.'
$ws.Range('D4').Value = 'No, The Indiana Oncology Society does not have direct influence on state or local policy. The organization primarily focuses on advancing cancer care and advocating for oncology professionals, rather than lobbying for policy changes.'
$ws.Range('E4').Value = 'No, IOS does not provide engagement opportunity with leadership. The focus of IOS is primarily on oncology education and networking among oncology professionals.'
$ws.Range('F4').Value = 'No, IOS does not provide support for clinical trial recruitment. IOS is an oncology society focused on education, networking, and advocacy for healthcare professionals in Indiana, rather than directly facilitating clinical trial recruitment.'
$ws.Range('G4').Value = 'No, IOS does not provide engagement opportunity with payors. IOS is a professional organization for oncologists and does not typically offer opportunities for engagement with payors.'
$ws.Range('H4').Value = 'Yes, 
The Indiana Oncology Society includes area experts in the field of oncology on its board to ensure comprehensive and informed decision-making.'
$ws.Range('I4').Value = 'No, IOS primarily focuses on education and advocacy for oncologists.'
$ws.Range('J4').Value = 'No, 
The Indiana Oncology Society may or may not include top therapeutic area experts on its board, but the structure of the society and its composition is not fully clear from the abbreviation "IOS" provided.'
$ws.Range('K4').Value = 'Midwest.'
$ws.Range('B5').Value = '132'
$ws.Range('C5').Value = 'No, the IOWA Oncology Society does not encompass community sites. It is focused on oncology professionals and research institutions in Iowa.'
$ws.Range('D5').Value = 'No, lack of public information or evidence of direct policy influence.'
$ws.Range('E5').Value = 'Yes, IOWA Oncology Society provides engagement opportunity with leadership. The society offers various opportunities for members to engage with leadership through conferences, workshops, and networking events.'
$ws.Range('F5').Value = 'Yes,  The Iowa Oncology Society provides support for clinical trial recruitment through education and resources for oncology professionals.'
$ws.Range('G5').Value = 'No, the IOWA Oncology Society does not provide engagement opportunities with payors. The organization''s focus is on oncology practice, research, and education, rather than payor relations.'
$ws.Range('H5').Value = 'Yes, the IOWA Oncology Society includes area experts on its board. This can be seen in the quality of information and guidance provided by the society.'
$ws.Range('I5').Value = 'Yes, The IOWA Oncology Society is involved in therapeutic research collaborations.'
$ws.Range('J5').Value = 'No, The IOWA Oncology Society does not have top therapeutic area experts on its board. This is evident from the make-up of the current board members who do not solely focus on a specific therapeutic area within oncology.'
$ws.Range('K5').Value = 'Midwest.'
$ws.Range('B6').Value = '450'
$ws.Range('C6').Value = 'Yes, MOASC does encompass community sites, as they are a regional association representing medical oncologists in Southern California, including those practicing in community settings.'
$ws.Range('D6').Value = 'No, MOASC does not have a direct influence on state or local policy as it primarily focuses on professional development and education for medical oncologists in Southern California.'
$ws.Range('E6').Value = 'Yes, MOASC provides engagement opportunity with leadership. MOASC offers various leadership positions and opportunities for members to become involved in shaping the organization''s direction and initiatives.'
$ws.Range('F6').Value = 'No, MOASC does not provide clinical trial recruitment support. They primarily focus on education, advocacy, and networking opportunities for medical oncologists in Southern California.'
$ws.Range('G6').Value = 'No, MOASC does not engage with payors. ,MOASC focuses on providing education and support for medical oncology professionals in Southern California.'
$ws.Range('H6').Value = 'Yes, MOASC does include area experts on its board. The organization is comprised of oncologists, pharmacists, nurses, and other healthcare professionals with expertise in medical oncology.'
$ws.Range('I6').Value = 'Yes, MOASC is involved in collaborative research efforts with various organizations in the field of oncology ,as evidenced by their participation in multiple clinical trials and research projects.'
$ws.Range('J6').Value = 'Yes, 
MOASC includes top therapeutic area experts on its board because its members consist of medical oncologists and other healthcare professionals specializing in oncology in Southern California.'
$ws.Range('K6').Value = 'Southern California'

# Restore default (unformatted) style on the numeric-text cells, keeping their
# values stored as text (matching the workbook author's original layout).
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = "Normal"
}